$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "60.762.00"
Set-TextValue $ws.Range("E2") "  -3.57%  "
Set-TextValue $ws.Range("D3") "2.907.92"
Set-TextValue $ws.Range("E3") "  -4.01%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "589.81"
Set-TextValue $ws.Range("E5") "  -0.80%  "
Set-TextValue $ws.Range("D6") "144.49"
Set-TextValue $ws.Range("E6") "  -5.76%  "
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("E8") "  -2.20%  "
Set-TextValue $ws.Range("D9") "2.907.31"
Set-TextValue $ws.Range("E9") "  -3.92%  "
Set-TextValue $ws.Range("E10") "  -3.73%  "
Set-TextValue $ws.Range("E11") "  -4.49%  "
Set-TextValue $ws.Range("D12") "0.444"
Set-TextValue $ws.Range("E13") "  -3.45%  "
Set-TextValue $ws.Range("D14") "33.40"
Set-TextValue $ws.Range("E14") "  -5.98%  "
Set-TextValue $ws.Range("D16") "3.386.96"
Set-TextValue $ws.Range("E16") "  -4.27%  "
Set-TextValue $ws.Range("D17") "60.721.60"
Set-TextValue $ws.Range("E17") "  -3.55%  "
Set-TextValue $ws.Range("E18") "  -5.72%  "
Set-TextValue $ws.Range("D19") "2.901.75"
Set-TextValue $ws.Range("E19") "  -4.25%  "
Set-TextValue $ws.Range("D20") "428.07"
Set-TextValue $ws.Range("E20") "  -5.19%  "
Set-TextValue $ws.Range("D21") "13.52"
Set-TextValue $ws.Range("E21") "  -5.28%  "
Set-TextValue $ws.Range("E22") "  -1.62%  "
Set-TextValue $ws.Range("D23") "7.07"
Set-TextValue $ws.Range("D24") "81.05"
Set-TextValue $ws.Range("E24") "  -2.63%  "
Set-TextValue $ws.Range("D25") "10.75"
Set-TextValue $ws.Range("E25") "  -6.44%  "
Set-TextValue $ws.Range("D26") "2.21"
Set-TextValue $ws.Range("E26") "  -5.26%  "
Set-TextValue $ws.Range("D27") "11.88"
Set-TextValue $ws.Range("E27") "  -4.16%  "
Set-TextValue $ws.Range("E28") "  +0.03%  "
Set-TextValue $ws.Range("D29") "2.23"
Set-TextValue $ws.Range("E29") "  -2.40%  "
Set-TextValue $ws.Range("D31") "2.62"
Set-TextValue $ws.Range("E31") "  -3.23%  "
Set-TextValue $ws.Range("D32") "7.10"
Set-TextValue $ws.Range("E32") "  -6.19%  "
Set-TextValue $ws.Range("D33") "26.42"
Set-TextValue $ws.Range("E33") "  -4.28%  "
Set-TextValue $ws.Range("E34") "  -4.25%  "
Set-TextValue $ws.Range("D35") "0.0₃0850"
Set-TextValue $ws.Range("E35") "  -1.36%  "
Set-TextValue $ws.Range("E36") "  -3.05%  "
Set-TextValue $ws.Range("E37") "  -5.04%  "
Set-TextValue $ws.Range("D38") "3.03"
Set-TextValue $ws.Range("E38") "  -3.27%  "
Set-TextValue $ws.Range("D39") "49.28"
Set-TextValue $ws.Range("E39") "  -2.38%  "
Set-TextValue $ws.Range("D40") "0.124"
Set-TextValue $ws.Range("E40") "  -6.41%  "
Set-TextValue $ws.Range("D41") "1.99"
Set-TextValue $ws.Range("E41") "  -5.32%  "
Set-TextValue $ws.Range("E42") "  -5.69%  "
Set-TextValue $ws.Range("D43") "0.293"
Set-TextValue $ws.Range("E43") "  -5.10%  "
Set-TextValue $ws.Range("D44") "40.60"
Set-TextValue $ws.Range("E44") "  -8.74%  "
Set-TextValue $ws.Range("E45") "  -3.00%  "
Set-TextValue $ws.Range("D46") "372.61"
Set-TextValue $ws.Range("E46") "  -4.96%  "
Set-TextValue $ws.Range("D47") "2.692.53"
Set-TextValue $ws.Range("E47") "  -1.06%  "
Set-TextValue $ws.Range("D48") "131.96"
Set-TextValue $ws.Range("E48") "  -0.72%  "
Set-TextValue $ws.Range("E49") "  +0.05%  "
Set-TextValue $ws.Range("D50") "24.04"
Set-TextValue $ws.Range("E50") "  -6.85%  "
Set-TextValue $ws.Range("E51") "  -3.10%  "
